# Auto commit at 2025-11-19  7:49:52.75
# Append two new daily rows (2025-11-18 data) to the charging-volume log on
# Sheet1: one row for the "四方坪站" station and one for the "高岭站" station,
# directly below the existing last row (157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 158: 四方坪站充电量(kw), date 2025-11-18 (serial 45979) ----
$ws.Range("A158").Value = 45979
$ws.Range("B158").Value = "四方坪站充电量(kw)"
$ws.Range("C158").Value = 528.7600000000001
$ws.Range("D158").Value = 1349.6229999999996
$ws.Range("E158").Value = 583.40499999999997
$ws.Range("F158").Value = 453.99200000000008
$ws.Range("G158").Value = 396.47
$ws.Range("H158").Value = 482.54700000000003
$ws.Range("I158").Value = 652.70900000000006
$ws.Range("J158").Value = 227.01200000000003
$ws.Range("K158").Value = 116.179
$ws.Range("L158").Value = 135.32999999999998
$ws.Range("M158").Value = 100.97000000000001
$ws.Range("N158").Value = 241.584
$ws.Range("O158").Value = 945.55299999999977
$ws.Range("P158").Value = 1418.0979999999995
$ws.Range("Q158").Value = 606.21
$ws.Range("R158").Value = 417.39100000000008
$ws.Range("S158").Value = 340.04499999999996
$ws.Range("T158").Value = 353.04399999999998
$ws.Range("U158").Value = 141.04000000000002
$ws.Range("V158").Value = 162.14400000000001
$ws.Range("W158").Value = 46.5
$ws.Range("X158").Value = 32.6
$ws.Range("Y158").Value = 159.76699999999997
$ws.Range("Z158").Value = 102.96000000000001

# ---- Row 159: 高岭站充电量(kw), date 2025-11-18 (serial 45979) ----
$ws.Range("A159").Value = 45979
$ws.Range("B159").Value = "高岭站充电量(kw)"
$ws.Range("C159").Value = 346.99200000000002
$ws.Range("D159").Value = 510.98899999999998
$ws.Range("E159").Value = 122.672
$ws.Range("F159").Value = 94.592999999999989
$ws.Range("G159").Value = 2.3519999999999999
$ws.Range("H159").Value = 5.7189999999999994
$ws.Range("I159").Value = 405.27799999999996
$ws.Range("J159").Value = 145.898
$ws.Range("K159").Value = 321.83200000000005
$ws.Range("L159").Value = 348.471
$ws.Range("M159").Value = 297.73699999999997
$ws.Range("N159").Value = 303.50899999999996
$ws.Range("O159").Value = 538.49899999999991
$ws.Range("P159").Value = 797.64200000000017
$ws.Range("Q159").Value = 391.89699999999999
$ws.Range("R159").Value = 410.05200000000008
$ws.Range("S159").Value = 297.37200000000001
$ws.Range("T159").Value = 70.186999999999983
$ws.Range("U159").Value = 59.524000000000001
$ws.Range("V159").Value = 87.137
$ws.Range("W159").Value = 60.516999999999996
$ws.Range("X159").Value = 125.52199999999999
$ws.Range("Y159").Value = 33.305
$ws.Range("Z159").Value = 15.484

# ---- Match the author's final cursor position (selection moved to E163) ----
$ws.Range("E163").Select() | Out-Null
